$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "gray" row's hex color codes (base and light shades)
$ws.Range("B3").Value = "#15343D"
$ws.Range("C3").Value = "#1F4D5B"

# Update the selection to reflect the active cell used when the file was saved
$ws.Range("C3").Select()
